$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates are Excel serial date numbers)
$data = @(
    @{ Row = 234; A = 44308; B = 2; C = 11; D = 335.8778625954199 },
    @{ Row = 235; A = 44309; B = 1; C = 7;  D = 213.7404580152672 },
    @{ Row = 236; A = 44310; B = 1; C = 5;  D = 152.6717557251908 },
    @{ Row = 237; A = 44311; B = 2; C = 6;  D = 183.206106870229 },
    @{ Row = 238; A = 44312; B = 3; C = 9;  D = 274.8091603053435 }
)

foreach ($item in $data) {
    $r = $item.Row

    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D

    # Copy the style from the row above (r-1) so formatting (date-style
    # border/alignment/number format on column A) matches existing rows.
    $ws.Cells.Item($r - 1, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
